# Edit script: updates NIT-9006348929.xlsx "Estado de Cuenta" sheet
# - Refreshes VALOR MORA total
# - Reduces worker count / period count
# - Replaces worker detail rows with new dataset (Estefany Acevedo Barrios + Richard Olivo Jimenez, 6 periods)
# - Removes now-obsolete rows for Alvaro Andres Mejia Acosta / Juan Rafael Reales Herrera
# - Adjusts column D width to fit new (shorter) longest name

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header summary values ---
$ws.Range("E11").Value = 208539      # VALOR MORA
$ws.Range("C13").Value = 2           # Cant. Trabajadores
$ws.Range("F13").Value = 7           # Cant. Periodos

# --- Remove the four middle detail rows (old rows 18-21) so the remaining ---
# --- rows keep correct borders/styles: 16,17 (regular top rows), 22-25 (regular),
# --- and 26 (bottom-border row) collapse down to rows 16-22.
$ws.Rows("18:21").Delete()

# --- Populate the 7 remaining detail rows with the new dataset ---

# Row 16: ESTEFANY ACEVEDO BARRIOS
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1143406693"
$ws.Range("D16").Value = "ESTEFANY ACEVEDO BARRIOS"
$ws.Range("E16").Value = "1811"
$ws.Range("F16").Value = 8333
$ws.Range("G16").Value = 781242

# Row 17: RICHARD OLIVO JIMENEZ - periodo 1901
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1143395492"
$ws.Range("D17").Value = "RICHARD OLIVO JIMENEZ"
$ws.Range("E17").Value = "1901"
$ws.Range("F17").Value = 31249
$ws.Range("G17").Value = 781242

# Row 18: RICHARD OLIVO JIMENEZ - periodo 1902
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1143395492"
$ws.Range("D18").Value = "RICHARD OLIVO JIMENEZ"
$ws.Range("E18").Value = "1902"
$ws.Range("F18").Value = 40000
$ws.Range("G18").Value = 781242

# Row 19: RICHARD OLIVO JIMENEZ - periodo 1903
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1143395492"
$ws.Range("D19").Value = "RICHARD OLIVO JIMENEZ"
$ws.Range("E19").Value = "1903"
$ws.Range("F19").Value = 40000
$ws.Range("G19").Value = 781242

# Row 20: RICHARD OLIVO JIMENEZ - periodo 1904
$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "1143395492"
$ws.Range("D20").Value = "RICHARD OLIVO JIMENEZ"
$ws.Range("E20").Value = "1904"
$ws.Range("F20").Value = 40000
$ws.Range("G20").Value = 781242

# Row 21: RICHARD OLIVO JIMENEZ - periodo 1905
$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "1143395492"
$ws.Range("D21").Value = "RICHARD OLIVO JIMENEZ"
$ws.Range("E21").Value = "1905"
$ws.Range("F21").Value = 31249
$ws.Range("G21").Value = 781242

# Row 22: RICHARD OLIVO JIMENEZ - periodo 1906
$ws.Range("B22").Value = "CC"
$ws.Range("C22").Value = "1143395492"
$ws.Range("D22").Value = "RICHARD OLIVO JIMENEZ"
$ws.Range("E22").Value = "1906"
$ws.Range("F22").Value = 17708
$ws.Range("G22").Value = 781242

# --- Column D is now narrower since the longest remaining name is shorter ---
# (best-fit width for the new, shorter longest name "ESTEFANY ACEVEDO BARRIOS")
$ws.Columns("D").ColumnWidth = 27.43
